# "merged jacobs comments into the structure"
#
# Fills in the peer-assessment grade ("Good") and the example-of-actions
# comment for Jacob in both criterion blocks (row 3 / Criterion 1 Online
# collaboration, and row 22 / Criterion 1 International Collaboration),
# and moves the sheet's active selection to B6 (scrolled back to the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")
$ws.Activate()

# --- Row 3: Jacob's grade + comment under "Criterion 1 Online collaboration"
$ws.Range("B3").Value = "Good"
$ws.Range("C3").Value = "1) Robert has been the main driven part of setting up, configuring and `nkeeping track of the GitHub platform and administrating Discord server. `n2) He has been very active on Discord, and been a clear communicator about`nupdates on GitHub and his contributions to a variety of project`nrelevant work. "
$ws.Range("C3").WrapText = $true

# --- Row 22: Jacob's grade + comment under "Criterion 1 International Collaboration"
$ws.Range("B22").Value = "Good"
$ws.Range("C22").Value = "1) Active collaboration on Discord and Git`n2) Good job on administrating both services "
$ws.Range("C22").WrapText = $true

# --- Move the selection / scroll position (was A7, now B6; topLeftCell A4 no longer pinned)
$ws.Range("B6").Select()
